$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 6: new dialogue line ---
$ws.Range("B6").Value = 219
$ws.Range("C6").Value = " Good to see you with us...[K]\nYour absence kept me from resting in peace…"
$ws.Range("D6").Value = " Рад снова вас видеть...[K] Я не\nнаходил себе покоя из-за вашего\nотсутствия..."
$ws.Range("E6").Value = " Ñàä òîïâà âàò âéäåóû...[K] Ÿ îå\nîàöïäéì òåáå ðïëïÿ éè-èà âàšåãï\nïóòôóòóâéÿ..."
$ws.Rows.Item(6).RowHeight = 38.4

# --- Give rows 5 and 6 together a thin top+bottom separator border ---
$rows56 = $ws.Range("A5:E6")
$rows56.Borders.Item(8).LineStyle = 1
$rows56.Borders.Item(8).Weight = 2
$rows56.Borders.Item(9).LineStyle = 1
$rows56.Borders.Item(9).Weight = 2

# --- Rows 7 & 8: two more dialogue lines (default, unbordered style) ---
$ws.Range("B7").Value = 178
$ws.Range("B8").Value = 181

$ws.Range("C7").Value = " My thoughts go with you…"
$ws.Range("C8").Value = " Please, put an end to the planet\'s\nparalysis."

$ws.Range("D7").Value = " Мои мысли только о вас..."
$ws.Range("D8").Value = " Прошу, положите конец\nпланетарному параличу."

$ws.Range("E7").Value = " Íïé íúòìé óïìûëï ï âàò..."
$ws.Range("E8").Value = " Ðñïšô, ðïìïçéóå ëïîåø\nðìàîåóàñîïíô ðàñàìéœô."

$ws.Rows.Item(8).RowHeight = 21.6

# --- Update the view: scroll so row 4 is at the top, select the last cell ---
$ws.Application.ActiveWindow.ScrollRow = 4
$ws.Range("E8").Select() | Out-Null

Write-Host "done"
